# Auto update stock data
# Updates the "Date_1" (A) and "EBITDA" (B) columns for each company's most
# recent row (2025/12/02 -> 2025/12/03), along with refreshed EBITDA figures.
#
# The date strings and EBITDA figures are stored as plain text in the
# workbook (inline strings), not as real dates/numbers. Assigning a
# numeric- or date-looking string straight to Range.Value would make Excel
# silently coerce it into a real number/date serial, so we briefly force
# the cell to Text format ("@") before the assignment, then put the cell
# style back to Normal so the rest of the formatting/layout is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  Date = "2025/12/03"; EBITDA = "5.07" },
    @{ Row = 8;  Date = "2025/12/03"; EBITDA = "7.68" },
    @{ Row = 14; Date = "2025/12/03"; EBITDA = "2.81" },
    @{ Row = 20; Date = "2025/12/03"; EBITDA = "12.33" },
    @{ Row = 26; Date = "2025/12/03"; EBITDA = "10.11" },
    @{ Row = 32; Date = "2025/12/03"; EBITDA = "25.99" },
    @{ Row = 38; Date = "2025/12/03" },
    @{ Row = 44; Date = "2025/12/03"; EBITDA = "11.10" },
    @{ Row = 50; Date = "2025/12/03"; EBITDA = "11.79" },
    @{ Row = 56; Date = "2025/12/03"; EBITDA = "31.87" },
    @{ Row = 62; Date = "2025/12/03"; EBITDA = "11.52" },
    @{ Row = 68; Date = "2025/12/03"; EBITDA = "12.04" },
    @{ Row = 74; Date = "2025/12/03"; EBITDA = "16.52" }
)

foreach ($u in $updates) {
    $r = $u.Row
    Set-TextValue $ws.Range("A$r") $u.Date
    if ($u.ContainsKey("EBITDA")) {
        Set-TextValue $ws.Range("B$r") $u.EBITDA
    }
}
